$d = $word.ActiveDocument

# --- Edit 1: "Isa Clemente" -> "Isabella" + " Clemente" (two separate runs) ---
$rng1 = $d.Content
$rng1.Find.Execute("Isa Clemente", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target1 = $d.Range($rng1.Start, $rng1.End)
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:r><w:rPr>' +
        '<w:rFonts w:ascii="Titillium Web" w:eastAsia="Titillium Web" w:hAnsi="Titillium Web" w:cs="Titillium Web"/>' +
        '<w:color w:val="666666"/><w:sz w:val="28"/><w:szCs w:val="28"/>' +
        '</w:rPr><w:t>Isabella</w:t></w:r>' +
        '<w:r><w:rPr>' +
        '<w:rFonts w:ascii="Titillium Web" w:eastAsia="Titillium Web" w:hAnsi="Titillium Web" w:cs="Titillium Web"/>' +
        '<w:color w:val="666666"/><w:sz w:val="28"/><w:szCs w:val="28"/>' +
        '</w:rPr><w:t xml:space="preserve"> Clemente</w:t></w:r>' +
        '</w:p>'
$target1.InsertXML($xml1)

# --- Edit 2: "Siopen - Altamura" -> "Siopen" (wrapped in proofErr spellStart/spellEnd) + " - Altamura" ---
$rng2 = $d.Content
$rng2.Find.Execute("Siopen - Altamura", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target2 = $d.Range($rng2.Start, $rng2.End)
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr>' +
        '<w:rFonts w:ascii="Titillium Web" w:eastAsia="Titillium Web" w:hAnsi="Titillium Web" w:cs="Titillium Web"/>' +
        '<w:color w:val="666666"/><w:sz w:val="28"/><w:szCs w:val="28"/>' +
        '</w:rPr><w:t>Siopen</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr>' +
        '<w:rFonts w:ascii="Titillium Web" w:eastAsia="Titillium Web" w:hAnsi="Titillium Web" w:cs="Titillium Web"/>' +
        '<w:color w:val="666666"/><w:sz w:val="28"/><w:szCs w:val="28"/>' +
        '</w:rPr><w:t xml:space="preserve"> - Altamura</w:t></w:r>' +
        '</w:p>'
$target2.InsertXML($xml2)
